$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (current "Tipo" column), shifting
# "Tipo" to column E, and making room for the new "MAE" column.
$ws.Range("D1").EntireColumn.Insert()

# Populate the new column D with header and value, matching the style
# used by the other header cells (copy style from existing header C1).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.1560271560555998

$wb.Save()
